$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143 (table-like data block "Chirimoya" / Vega Modelo de Temuco),
# shifting the existing rows 143-222 down to 144-223.
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with the new record's data.
$ws.Range("A143").Value2 = 10
$ws.Range("B143").Value2 = "Vega Modelo de Temuco"
$ws.Range("C143").Value2 = "La Araucanía"
$ws.Range("D143").Value2 = 45202
$ws.Range("E143").Value2 = 9
$ws.Range("F143").Value2 = "Fruta"
$ws.Range("G143").Value2 = 100107
$ws.Range("H143").Value2 = "Otros"
$ws.Range("I143").Value2 = 100107002
$ws.Range("J143").Value2 = "Chirimoya"
$ws.Range("K143").Value2 = "Cultivar IV Región"
$ws.Range("L143").Value2 = "Primera"
$ws.Range("M143").Value2 = 100
$ws.Range("N143").Value2 = 2500
$ws.Range("O143").Value2 = 2500
$ws.Range("P143").Value2 = 2500
$ws.Range("Q143").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R143").Value2 = "Provincia del Elquí"
$ws.Range("S143").Value2 = 2500
$ws.Range("T143").Value2 = 1
